# Applies the diff:
#  1. Rearranges (swaps/rotates) the F:V ("match") data among certain existing
#     rows (A:E - Indice/pais/torneio/temporada/data_partida - stay untouched).
#  2. Appends six brand-new match rows (73-78) at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Get-RowData {
    param($ws, [int]$r)
    $data = @{}
    for ($c = 6; $c -le 22; $c++) {
        # columns F(6) .. V(22)
        $data[$c] = $ws.Cells.Item($r, $c).Value()
    }
    return $data
}

function Set-RowData {
    param($ws, [int]$r, $data)
    for ($c = 6; $c -le 22; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c]
    }
}

# ---------------------------------------------------------------------------
# 1) Re-arrange existing rows: new row <dest> = old row <src> (columns F:V)
# ---------------------------------------------------------------------------

$rowMap = @{
    7 = 9; 9 = 7;
    13 = 14; 14 = 13;
    16 = 17; 17 = 16;
    20 = 21; 21 = 20;
    29 = 30; 30 = 31; 31 = 29;
    40 = 42; 41 = 40; 42 = 41;
    46 = 50; 47 = 49; 49 = 47; 50 = 46;
    52 = 55; 53 = 54; 54 = 53; 55 = 52;
    69 = 70; 70 = 71; 71 = 69;
}

# snapshot every source row's data *before* overwriting anything, since many
# rows are both a source and a destination.
$snapshots = @{}
foreach ($src in $rowMap.Values) {
    if (-not $snapshots.ContainsKey($src)) {
        $snapshots[$src] = Get-RowData $ws $src
    }
}

foreach ($dest in $rowMap.Keys) {
    $src = $rowMap[$dest]
    Set-RowData $ws $dest $snapshots[$src]
}

# ---------------------------------------------------------------------------
# 2) Append new rows 73-78
# ---------------------------------------------------------------------------

# Copy formatting (styles) from the last existing data row (72) down to the
# six new rows so they match the rest of the table (bold/border index on A,
# datetime format on E, etc.).
$ws.Range("A72:V72").Copy()
$ws.Range("A73:V78").PasteSpecial(-4122)  # xlPasteFormats

$newRows = @(
    @{ Row = 73; A = 72; E = 45242.5625;          F = "PAE Egaleo"; G = 1; H = "Giouchtas";        I = 0;
       J = 1.83; K = "11/11/2023 01:42"; L = 2.07; M = "12/11/2023 13:26";
       N = 3.18; O = "11/11/2023 01:42"; P = 3.11; Q = "12/11/2023 13:26";
       R = 4.11; S = "11/11/2023 01:42"; T = 3.75; U = "12/11/2023 13:26";
       V = "https://www.betexplorer.com/football/greece/super-league-2/pae-egaleo-giouchtas/lCDgX94L/" },

    @{ Row = 74; A = 73; E = 45242.58333333334;   F = "Panachaiki"; G = 1; H = "Ionikos";           I = 1;
       J = 4.2;  K = "11/11/2023 02:12"; L = 6.83; M = "12/11/2023 13:58";
       N = 3.26; O = "11/11/2023 02:12"; P = 3.77; Q = "12/11/2023 13:58";
       R = 1.75; S = "11/11/2023 02:12"; T = 1.52; U = "12/11/2023 13:58";
       V = "https://www.betexplorer.com/football/greece/super-league-2/panachaiki-ionikos/2JgUQ7sk/" },

    @{ Row = 75; A = 74; E = 45242.58333333334;   F = "Levadiakos"; G = 2; H = "Aiolikos";          I = 0;
       J = 1.24; K = "11/11/2023 04:43"; L = 1.27; M = "12/11/2023 11:59";
       N = 4.72; O = "11/11/2023 04:43"; P = 4.98; Q = "12/11/2023 12:05";
       R = 10.75; S = "11/11/2023 04:43"; T = 12.44; U = "12/11/2023 11:59";
       V = "https://www.betexplorer.com/football/greece/super-league-2/levadiakos-aiolikos-fc/URrrJbdC/" },

    @{ Row = 76; A = 75; E = 45242.58333333334;   F = "Ilioupoli"; G = 1; H = "Athens Kallithea";   I = 2;
       J = 4.68; K = "11/11/2023 02:12"; L = 4.56; M = "12/11/2023 13:55";
       N = 3.36; O = "11/11/2023 02:12"; P = 3.13; Q = "12/11/2023 13:55";
       R = 1.65; S = "11/11/2023 02:12"; T = 1.88; U = "12/11/2023 13:55";
       V = "https://www.betexplorer.com/football/greece/super-league-2/ilioupoli-athens-kallithea/GMcYPRde/" },

    @{ Row = 77; A = 76; E = 45242.58333333334;   F = "PAE Chania"; G = 2; H = "Panathinaikos B";   I = 0;
       J = 1.37; K = "11/11/2023 02:12"; L = 1.32; M = "12/11/2023 12:05";
       N = 4.28; O = "11/11/2023 02:12"; P = 4.79; Q = "12/11/2023 12:07";
       R = 6.55; S = "11/11/2023 02:12"; T = 9.890000000000001; U = "12/11/2023 12:07";
       V = "https://www.betexplorer.com/football/greece/super-league-2/pae-chania-panathinaikos/8dfQRmSr/" },

    @{ Row = 78; A = 77; E = 45242.58333333334;   F = "Niki Volos"; G = 3; H = "Iraklis 1908";      I = 1;
       J = 1.87; K = "11/11/2023 02:12"; L = 2.02; M = "12/11/2023 12:05";
       N = 3.03; O = "11/11/2023 02:12"; P = 3;    Q = "12/11/2023 12:12";
       R = 3.98; S = "11/11/2023 02:12"; T = 4.13; U = "12/11/2023 12:05";
       V = "https://www.betexplorer.com/football/greece/super-league-2/niki-volos-iraklis-fc/jypvKvs6/" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A            # A - Indice
    $ws.Cells.Item($r, 2).Value = "greece"           # B - pais
    $ws.Cells.Item($r, 3).Value = "super-league-2"   # C - torneio
    $ws.Cells.Item($r, 4).Value = "2023-2024"        # D - temporada
    $ws.Cells.Item($r, 5).Value = $row.E             # E - data_partida
    $ws.Cells.Item($r, 6).Value = $row.F             # F - home
    $ws.Cells.Item($r, 7).Value = $row.G             # G - home_ft_gols
    $ws.Cells.Item($r, 8).Value = $row.H             # H - away
    $ws.Cells.Item($r, 9).Value = $row.I             # I - away_ft_gols
    $ws.Cells.Item($r, 10).Value = $row.J            # J - home_opening_odds
    $ws.Cells.Item($r, 11).Value = $row.K            # K - home_opening_data_hora
    $ws.Cells.Item($r, 12).Value = $row.L            # L - home_closing_odds
    $ws.Cells.Item($r, 13).Value = $row.M            # M - home_closing_data_hora
    $ws.Cells.Item($r, 14).Value = $row.N            # N - draw_opening_odds
    $ws.Cells.Item($r, 15).Value = $row.O            # O - draw_opening_data_hora
    $ws.Cells.Item($r, 16).Value = $row.P            # P - draw_closing_odds
    $ws.Cells.Item($r, 17).Value = $row.Q            # Q - draw_closing_data_hora
    $ws.Cells.Item($r, 18).Value = $row.R            # R - away_opening_odds
    $ws.Cells.Item($r, 19).Value = $row.S            # S - away_opening_data_hora
    $ws.Cells.Item($r, 20).Value = $row.T            # T - away_closing_odds
    $ws.Cells.Item($r, 21).Value = $row.U            # U - away_closing_data_hora
    $ws.Cells.Item($r, 22).Value = $row.V            # V - url_partida
}
